# daily auto push: 2026-01-30 22:44 UTC
#
# Inserts one new record row ("2026/01/31", "土", 3, 201) right before the
# existing 2026/12/29 block (old row 743), shifting every row from the old
# 743 down through the old 784 (the last row) down by one — they become
# rows 744..785. The sheet's used-range dimension grows from D784 to D785.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything at/after row 743 down by one row.
$ws.Rows(743).Insert()

# Column A holds "yyyy/mm/dd" looking text that Excel's COM layer would
# otherwise auto-convert to a real date serial the moment a plain .Value
# (or .Formula/.FormulaR1C1) assignment hits a "General" formatted cell.
# Force the cell into Text mode first so the literal string sticks, then
# drop the explicit number format again so the cell is left unstyled —
# exactly like every other date cell in this column.
$ws.Range("A743").NumberFormat = "@"
$ws.Range("A743").Value = "2026/01/31"
$ws.Range("A743").ClearFormats()

$ws.Range("B743").Value = "土"
$ws.Range("C743").Value = 3
$ws.Range("D743").Value = 201
